$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "68.990.22"
    "E2" = "  -0.12%  "
    "D3" = "3.767.88"
    "E3" = "  -1.65%  "
    "D4" = "0.999"
    "E4" = "  -0.36%  "
    "D5" = "630.35"
    "E5" = "  -0.87%  "
    "D6" = "165.88"
    "E6" = "  +0.28%  "
    "D7" = "3.766.55"
    "E7" = "  -1.63%  "
    "E8" = "  +0.05%  "
    "E9" = "  -0.13%  "
    "E10" = "  -2.11%  "
    "E11" = "  +0.57%  "
    "D12" = "6.75"
    "E12" = "  +0.98%  "
    "E13" = "  -4.36%  "
    "D14" = "34.86"
    "E14" = "  -3.08%  "
    "D15" = "4.401.16"
    "E15" = "  -1.63%  "
    "D16" = "3.765.11"
    "E16" = "  -1.92%  "
    "D17" = "68.958.64"
    "D18" = "17.66"
    "E18" = "  -3.33%  "
    "E19" = "  +0.06%  "
    "D20" = "7.03"
    "E20" = "  -1.85%  "
    "D21" = "462.24"
    "E21" = "  -1.02%  "
    "D22" = "9.50"
    "E22" = "  -2.68%  "
    "D23" = "0.704"
    "E23" = "  -0.75%  "
    "E24" = "  -5.46%  "
    "D25" = "82.09"
    "E25" = "  -2.19%  "
    "E26" = "  +0.45%  "
    "D27" = "2.13"
    "E27" = "  -1.43%  "
    "D28" = "10.13"
    "E28" = "  +0.49%  "
    "E29" = "  -0.04%  "
    "D30" = "3.918.21"
    "E30" = "  -1.58%  "
    "E31" = "  +0.83%  "
    "E32" = "  -0.71%  "
    "D33" = "7.06"
    "E33" = "  -3.67%  "
    "D34" = "28.33"
    "E34" = "  -2.99%  "
    "D35" = "0.175"
    "E35" = "  +17.99%  "
    "E36" = "  +0.10%  "
    "D37" = "3.719.71"
    "E37" = "  -1.55%  "
    "D38" = "8.90"
    "E38" = "  -2.53%  "
    "E39" = "  -1.49%  "
    "D40" = "3.29"
    "E40" = "  +0.22%  "
    "E41" = "  -2.71%  "
    "D42" = "0.999"
    "E42" = "  -0.16%  "
    "E43" = "  -2.15%  "
    "D45" = "156.67"
    "E45" = "  -0.10%  "
    "E46" = "  +4.10%  "
    "D47" = "1.42"
    "E47" = "  -0.71%  "
    "D48" = "46.94"
    "E48" = "  +0.17%  "
    "D49" = "42.87"
    "E49" = "  -0.08%  "
    "D50" = "0.294"
    "E50" = "  -2.45%  "
    "D51" = "8.35"
    "E51" = "  -0.87%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text interpretation so purely-numeric-looking strings
    # (e.g. "0.999", "9.50") are not auto-converted to numbers,
    # and non-numeric strings (e.g. "68.990.22") stay text as before.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Reset the style back to the workbook default so no stray
    # "Text" number-format style is left attached to the cell.
    $cell.Style = "Normal"
}
